# Update gh-pages output data (广州-漫展信息.xlsx)
# Applies refreshed "want to go" counters / sold-out labels scraped from
# bilibili show listings, and replaces the outdated duplicate
# "Look Look动漫嘉年华" row with the next two upcoming events, dropping the
# now-redundant last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1Numbers = @{
    "F2"  = 37122
    "F4"  = 630
    "F5"  = 761
    "F6"  = 474
    "F7"  = 159
    "F11" = 683
    "F12" = 520
    "F16" = 464
    "F17" = 438
    "F18" = 1152
    "F20" = 800
    "F21" = 2478
    "F22" = 976
    "F23" = 544
    "F24" = 98
    "F25" = 1149
    "F27" = 750
}
foreach ($cellRef in $ws1Numbers.Keys) {
    $ws1.Range($cellRef).Value = $ws1Numbers[$cellRef]
}
$ws1.Range("G2").Value = "已售罄"

# Row 28 used to be a duplicate of row 27 (Look Look动漫嘉年华). It now
# takes on the event that used to live in row 29 (第五届AP动漫嘉年华),
# with its "want to go" counter refreshed.
$ws1.Cells.Item(28, 3).Value = "广州·第五届AP动漫嘉年华"
$ws1.Cells.Item(28, 4).Value = "西环路1号 广州岭南会展中心"
$ws1.Cells.Item(28, 5).Value = "2024.06.01 10:00-06.01 17:00"
$ws1.Cells.Item(28, 6).Value = 47
$ws1.Cells.Item(28, 7).Value = 55
$ws1.Cells.Item(28, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws1.Cells.Item(28, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Row 29 takes on the event that used to live in row 30 (622排球少年only),
# with its "want to go" counter refreshed.
# (Force text so the date-shaped string isn't auto-converted to a date
# serial, matching how the rest of column B is stored.)
$ws1.Cells.Item(29, 2).NumberFormat = "@"
$ws1.Cells.Item(29, 2).Value = "2024-06-22"
$ws1.Cells.Item(29, 3).Value = "广州·622排球少年only"
$ws1.Cells.Item(29, 4).Value = "岭南购物城内 广州OMG网红街"
$ws1.Cells.Item(29, 5).Value = "2024.06.22 10:00-06.22 17:30"
$ws1.Cells.Item(29, 6).Value = 1141
$ws1.Cells.Item(29, 7).Value = 68
$ws1.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws1.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

# The old row 30 (now fully absorbed into row 29 above) is removed.
$ws1.Rows.Item(30).Delete()

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2Numbers = @{
    "F3" = 373
    "F9" = 8
}
foreach ($cellRef in $ws2Numbers.Keys) {
    $ws2.Range($cellRef).Value = $ws2Numbers[$cellRef]
}

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 614

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - union of the three sheets above
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4Numbers = @{
    "F2"  = 614
    "F3"  = 37122
    "F5"  = 630
    "F6"  = 761
    "F7"  = 474
    "F9"  = 159
    "F11" = 373
    "F16" = 683
    "F17" = 520
    "F22" = 8
    "F26" = 464
    "F27" = 438
    "F28" = 1152
    "F30" = 800
    "F31" = 2478
    "F32" = 976
    "F33" = 544
    "F34" = 98
    "F35" = 1149
    "F38" = 750
}
foreach ($cellRef in $ws4Numbers.Keys) {
    $ws4.Range($cellRef).Value = $ws4Numbers[$cellRef]
}
$ws4.Range("G3").Value = "已售罄"

# Row 39 used to be a duplicate of row 38 (Look Look动漫嘉年华). It now
# takes on the event that used to live in row 40 (第五届AP动漫嘉年华),
# with its "want to go" counter refreshed.
$ws4.Cells.Item(39, 3).Value = "广州·第五届AP动漫嘉年华"
$ws4.Cells.Item(39, 4).Value = "西环路1号 广州岭南会展中心"
$ws4.Cells.Item(39, 5).Value = "2024.06.01 10:00-06.01 17:00"
$ws4.Cells.Item(39, 6).Value = 47
$ws4.Cells.Item(39, 7).Value = 55
$ws4.Cells.Item(39, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws4.Cells.Item(39, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Row 40 takes on the event that used to live in row 41 (622排球少年only),
# with its "want to go" counter refreshed.
# (Force text so the date-shaped string isn't auto-converted to a date
# serial, matching how the rest of column B is stored.)
$ws4.Cells.Item(40, 2).NumberFormat = "@"
$ws4.Cells.Item(40, 2).Value = "2024-06-22"
$ws4.Cells.Item(40, 3).Value = "广州·622排球少年only"
$ws4.Cells.Item(40, 4).Value = "岭南购物城内 广州OMG网红街"
$ws4.Cells.Item(40, 5).Value = "2024.06.22 10:00-06.22 17:30"
$ws4.Cells.Item(40, 6).Value = 1141
$ws4.Cells.Item(40, 7).Value = 68
$ws4.Cells.Item(40, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82974"
$ws4.Cells.Item(40, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/WMlOXSZn1710748067155.jpeg"

# The old row 41 (now fully absorbed into row 40 above) is removed.
$ws4.Rows.Item(41).Delete()
